$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 0.06515966666666667
$ws.Range("M2").Value = 255.0443116666667
$ws.Range("N2").Value = 765.132935
$ws.Range("O2").Value = 0.863617428561108
$ws.Range("P2").Value = 0.8636174285611079
$ws.Range("Q2").Value = 16.61860233342945
$ws.Range("R2").Value = 149.567421000865
$ws.Range("S2").Value = 0.863617428561108
$ws.Range("T2").Value = 0.8636174285611079

# Row 3
$ws.Range("G3").Value = 0.06515966666666667
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("O3").Value = 0.003031431940796009
$ws.Range("P3").Value = 0.003031431940796009
$ws.Range("Q3").Value = 0.0583338875048889
$ws.Range("R3").Value = 0.5250049875440002
$ws.Range("S3").Value = 0.003031431940796009
$ws.Range("T3").Value = 0.003031431940796009

# Row 4
$ws.Range("G4").Value = 0.06515966666666667
$ws.Range("M4").Value = 7.050555333333333
$ws.Range("N4").Value = 21.151666
$ws.Range("O4").Value = 0.02387421396349043
$ws.Range("P4").Value = 0.02387421396349043
$ws.Range("Q4").Value = 0.4594118353348889
$ws.Range("R4").Value = 4.134706518014
$ws.Range("S4").Value = 0.02387421396349043
$ws.Range("T4").Value = 0.02387421396349043

# Row 5
$ws.Range("G5").Value = 0.06515966666666667
$ws.Range("M5").Value = 32.33082866666666
$ws.Range("N5").Value = 96.99248599999999
$ws.Range("O5").Value = 0.1094769255346056
$ws.Range("P5").Value = 0.1094769255346056
$ws.Range("Q5").Value = 2.106666018977111
$ws.Range("R5").Value = 18.959994170794
$ws.Range("S5").Value = 0.1094769255346056
$ws.Range("T5").Value = 0.1094769255346056
